# L9737_底稿_催收戶餘額明細.xlsx - add two new header columns to the
# "yyymmdd" sheet's header row: 評估淨值 (N1) and 貸放成數 (O1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yyymmdd")

# Match the on-screen zoom used when the column was added.
$excel.ActiveWindow.Zoom = 55

# Give the two new header cells the same look (bold/centred header style)
# as the rest of row 1 by copying the formatting from the last existing
# header cell (M1) before filling in the new labels.
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1:O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("N1").Value = "評估淨值"
$ws.Range("O1").Value = "貸放成數"

# Leave the selection on the newly added column, as in the saved file.
$ws.Range("O2").Select() | Out-Null
